$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a text-like cell value while forcing text storage
# (avoids Excel auto-parsing strings like "1.005" as numbers)
function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "21.237.97"
Set-TextCell $ws.Range("E2") "  +4.16%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.542.12"
Set-TextCell $ws.Range("E3") "  +5.54%  "

# Row 4
Set-TextCell $ws.Range("D4") "1.005"
Set-TextCell $ws.Range("E4") "  -0.42%  "

# Row 5
Set-TextCell $ws.Range("D5") "0.9591"
Set-TextCell $ws.Range("E5") "  +0.96%  "

# Row 6
Set-TextCell $ws.Range("D6") "282.24"
Set-TextCell $ws.Range("E6") "  +2.64%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.3628"
Set-TextCell $ws.Range("E7") "  -0.56%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.3193"
Set-TextCell $ws.Range("E8") "  +3.92%  "

# Row 9
Set-TextCell $ws.Range("D9") "40.96"
Set-TextCell $ws.Range("E9") "  +3.42%  "

# Row 10
Set-TextCell $ws.Range("D10") "1.098"
Set-TextCell $ws.Range("E10") "  +6.00%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.06811"
Set-TextCell $ws.Range("E11") "  +3.61%  "

# Row 12
Set-TextCell $ws.Range("D12") "1.000"
Set-TextCell $ws.Range("E12") "  -0.08%  "

# Row 13
Set-TextCell $ws.Range("D13") "5.678"
Set-TextCell $ws.Range("E13") "  +4.74%  "

# Row 14
Set-TextCell $ws.Range("D14") "18.80"
Set-TextCell $ws.Range("E14") "  +5.42%  "

# Row 15
Set-TextCell $ws.Range("D15") "6.361"
Set-TextCell $ws.Range("E15") "  +3.87%  "

# Row 16
Set-TextCell $ws.Range("D16") "0.00001051"
Set-TextCell $ws.Range("E16") "  +2.63%  "

# Row 17
Set-TextCell $ws.Range("D17") "0.9593"
Set-TextCell $ws.Range("E17") "  -0.86%  "

# Row 18
Set-TextCell $ws.Range("D18") "1.532.61"
Set-TextCell $ws.Range("E18") "  +4.76%  "

# Row 19
Set-TextCell $ws.Range("D19") "0.06059"
Set-TextCell $ws.Range("E19") "  +4.60%  "

# Row 20
Set-TextCell $ws.Range("D20") "72.43"
Set-TextCell $ws.Range("E20") "  +4.06%  "

# Row 21
Set-TextCell $ws.Range("D21") "5.693"
Set-TextCell $ws.Range("E21") "  +4.74%  "

# Row 22
Set-TextCell $ws.Range("E22") "  +4.45%  "

# Row 23
Set-TextCell $ws.Range("D23") "11.37"
Set-TextCell $ws.Range("E23") "  +4.59%  "

# Row 24
Set-TextCell $ws.Range("D24") "2.306"
Set-TextCell $ws.Range("E24") "  +2.77%  "

# Row 25
Set-TextCell $ws.Range("D25") "21.314.16"
Set-TextCell $ws.Range("E25") "  +4.32%  "

# Row 26
Set-TextCell $ws.Range("D26") "148.09"
Set-TextCell $ws.Range("E26") "  +4.48%  "

# Row 27
Set-TextCell $ws.Range("D27") "2.218"
Set-TextCell $ws.Range("E27") "  +6.73%  "

# Row 28
Set-TextCell $ws.Range("D28") "17.75"
Set-TextCell $ws.Range("E28") "  +3.57%  "

# Row 29
Set-TextCell $ws.Range("D29") "1.703.18"
Set-TextCell $ws.Range("E29") "  +5.37%  "

# Row 30
Set-TextCell $ws.Range("D30") "118.21"
Set-TextCell $ws.Range("E30") "  +5.58%  "

# Row 31
Set-TextCell $ws.Range("D31") "4.028"
Set-TextCell $ws.Range("E31") "  +5.25%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws.Range("D32") "5.234"
Set-TextCell $ws.Range("E32") "  +7.34%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextCell $ws.Range("D33") "0.8564"
Set-TextCell $ws.Range("E33") "  +8.38%  "

# Row 34
Set-TextCell $ws.Range("D34") "0.08015"
Set-TextCell $ws.Range("E34") "  +1.92%  "

# Row 35
Set-TextCell $ws.Range("D35") "1.505"
Set-TextCell $ws.Range("E35") "  -0.80%  "

# Row 36
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws.Range("D36") "4.990"
Set-TextCell $ws.Range("E36") "  +6.84%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Range("D37") "1.214"
Set-TextCell $ws.Range("E37") "  +6.43%  "

# Row 38
Set-TextCell $ws.Range("D38") "0.05877"
Set-TextCell $ws.Range("E38") "  +3.16%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.02108"
Set-TextCell $ws.Range("E39") "  +3.90%  "

# Row 40
Set-TextCell $ws.Range("D40") "10.74"
Set-TextCell $ws.Range("E40") "  +4.05%  "

# Row 41
Set-TextCell $ws.Range("D41") "7.745"
Set-TextCell $ws.Range("E41") "  +2.96%  "

# Row 42
Set-TextCell $ws.Range("D42") "0.1919"
Set-TextCell $ws.Range("E42") "  +3.65%  "

# Row 43
Set-TextCell $ws.Range("D43") "0.9593"
Set-TextCell $ws.Range("E43") "  +0.27%  "

# Row 44
Set-TextCell $ws.Range("D44") "0.5470"
Set-TextCell $ws.Range("E44") "  +4.03%  "

# Row 45
Set-TextCell $ws.Range("D45") "12.48"
Set-TextCell $ws.Range("E45") "  +4.49%  "

# Row 46
Set-TextCell $ws.Range("D46") "3.579"
Set-TextCell $ws.Range("E46") "  +2.62%  "

# Row 47
Set-TextCell $ws.Range("D47") "0.5461"
Set-TextCell $ws.Range("E47") "  +6.43%  "

# Row 48
Set-TextCell $ws.Range("D48") "122.04"
Set-TextCell $ws.Range("E48") "  +4.11%  "

# Row 49
Set-TextCell $ws.Range("D49") "1.875"
Set-TextCell $ws.Range("E49") "  +7.17%  "

# Row 50
Set-TextCell $ws.Range("D50") "0.06634"
Set-TextCell $ws.Range("E50") "  +3.51%  "

# Row 51
Set-TextCell $ws.Range("D51") "70.18"
Set-TextCell $ws.Range("E51") "  +6.09%  "
